$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Type" header column
$ws.Range("B1").Value = "Type"

# Counts for the existing numeric asset rows (2-7)
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1

# New asset-type rows appended at the bottom, in shared-string insertion order
$ws.Range("A10").Value = "0000000"
$ws.Range("B10").Value = 1

$ws.Range("A11").Value = "1111111"
$ws.Range("B11").Value = 1

# Insert a new row at position 8, pushing the old "1578465" / "00456789" rows
# and the two rows just appended down by one (preserves their value types).
$ws.Rows("8:8").Insert()

# Fill the freshly inserted row with a new text asset id + its count
$ws.Range("A8").Value = "2222222"
$ws.Range("B8").Value = 2

# Counts for the rows that shifted down but didn't have column B yet
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 1

$ws.Range("D9").Select()
